$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-10-04 12:39:34"
}
